$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = -21.99510000000001
$ws.Range("A21").Value = -19.78979999999998
$ws.Range("A23").Value = -20.36379999999997
$ws.Range("A25").Value = -21.60319999999999
$ws.Range("C27").Value = -12.615
$ws.Range("C31").Value = -13.004
$ws.Range("C39").Value = -12.75230000000001
$ws.Range("C48").Value = -11.66339999999999
$ws.Range("C51").Value = -11.9042
$ws.Range("C52").Value = -10.98589999999999
$ws.Range("A53").Value = -21.8668
$ws.Range("C55").Value = -13.78169999999999
$ws.Range("C56").Value = -12.5209
$ws.Range("A57").Value = -22.02659999999999
$ws.Range("C57").Value = -13.05209999999999
$ws.Range("A59").Value = -22.2561
$ws.Range("A69").Value = -21.62830000000001
$ws.Range("C73").Value = -12.38300000000001
$ws.Range("A79").Value = -20.43450000000001
$ws.Range("A83").Value = -21.9454
$ws.Range("C89").Value = -10.33890000000001
$ws.Range("C90").Value = -12.5555
$ws.Range("A93").Value = -21.3215
